# Update "想去人数" (F) and "最低票价" (G) figures across the three sheets
# that list the same 2024 Hangzhou comic-con events (展览 / 演出 / 全部类型),
# matching the gh-pages data refresh commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" ---
$ws1.Range("F2").Value = 5894
$ws1.Range("F3").Value = 556
$ws1.Range("G3").Value = 148
$ws1.Range("F4").Value = 1086
$ws1.Range("F5").Value = 1049
$ws1.Range("F7").Value = 87
$ws1.Range("F10").Value = 59
$ws1.Range("F11").Value = 29
$ws1.Range("F13").Value = 2039
$ws1.Range("F14").Value = 1522
$ws1.Range("F15").Value = 1116
$ws1.Range("F16").Value = 306
$ws1.Range("G16").Value = "已售罄"
$ws1.Range("F18").Value = 428
$ws1.Range("F19").Value = 652
$ws1.Range("F20").Value = 235
$ws1.Range("F24").Value = 3672
$ws1.Range("F25").Value = 199
$ws1.Range("F26").Value = 135
$ws1.Range("F28").Value = 167
$ws1.Range("F30").Value = 523
$ws1.Range("F36").Value = 847
$ws1.Range("F39").Value = 87
$ws1.Range("F40").Value = 89

# --- Sheet "演出" ---
$ws2.Range("G2").Value = 68
$ws2.Range("F4").Value = 754
$ws2.Range("F6").Value = 413

# --- Sheet "全部类型" ---
$ws4.Range("G2").Value = 68
$ws4.Range("F3").Value = 5894
$ws4.Range("F4").Value = 556
$ws4.Range("G4").Value = 148
$ws4.Range("F5").Value = 1086
$ws4.Range("F7").Value = 754
$ws4.Range("F8").Value = 1049
$ws4.Range("F11").Value = 413
$ws4.Range("F12").Value = 87
$ws4.Range("F15").Value = 59
$ws4.Range("F16").Value = 29
$ws4.Range("F19").Value = 2039
$ws4.Range("F20").Value = 1522
$ws4.Range("F21").Value = 1116
$ws4.Range("F22").Value = 306
$ws4.Range("G22").Value = "已售罄"
$ws4.Range("F24").Value = 428
$ws4.Range("F26").Value = 652
$ws4.Range("F27").Value = 235
$ws4.Range("F30").Value = 3672
$ws4.Range("F31").Value = 199
$ws4.Range("F32").Value = 135
$ws4.Range("F34").Value = 167
$ws4.Range("F36").Value = 523
$ws4.Range("F42").Value = 847
$ws4.Range("F45").Value = 87
$ws4.Range("F46").Value = 89
